# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest generated data (matching commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F
$updates = @{
    2  = 8446
    3  = 8084
    4  = 138
    8  = 138
    10 = 187
    12 = 729
    13 = 179
    14 = 2182
    16 = 68
    20 = 87
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
